# Auto-generated update of leve-profit market data cells across all class sheets.
# Source: scheduled market-data refresh (chore: update Sheets via scheduled runner)
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC (47 cell updates) ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 97.88
$ws.Range("I15").Value = 97.88
$ws.Range("K15").Value = 293.64
$ws.Range("M15").Value = -124.64
$ws.Range("H17").Value = 595294.4399999999
$ws.Range("J17").Value = 595294.4399999999
$ws.Range("L17").Value = 1785883.32
$ws.Range("N17").Value = -1786219.32
$ws.Range("H33").Value = 499
$ws.Range("I33").Value = 498.25
$ws.Range("K33").Value = 498.25
$ws.Range("M33").Value = -269.25
$ws.Range("H43").Value = 1198.8
$ws.Range("I43").Value = 833.3333
$ws.Range("J43").Value = 1355.4286
$ws.Range("K43").Value = 833.3333
$ws.Range("L43").Value = 1355.4286
$ws.Range("M43").Value = -764.3333
$ws.Range("N43").Value = -1493.4286
$ws.Range("H62").Value = 74901.14
$ws.Range("I62").Value = 169600.83
$ws.Range("J62").Value = 3876.375
$ws.Range("K62").Value = 169600.83
$ws.Range("L62").Value = 3876.375
$ws.Range("M62").Value = -168976.83
$ws.Range("N62").Value = -5124.375
$ws.Range("H65").Value = 74901.14
$ws.Range("I65").Value = 169600.83
$ws.Range("J65").Value = 3876.375
$ws.Range("K65").Value = 848004.1499999999
$ws.Range("L65").Value = 19381.875
$ws.Range("M65").Value = -844884.1499999999
$ws.Range("N65").Value = -25621.875
$ws.Range("H135").Value = 34885132
$ws.Range("I135").Value = 14286905
$ws.Range("J135").Value = 125002380
$ws.Range("K135").Value = 128582145
$ws.Range("L135").Value = 1125021420
$ws.Range("M135").Value = -128579610
$ws.Range("N135").Value = -1125026490
$ws.Range("H138").Value = 3574.2
$ws.Range("I138").Value = 2074.5417
$ws.Range("J138").Value = 4164.2295
$ws.Range("K138").Value = 6223.625100000001
$ws.Range("L138").Value = 12492.6885
$ws.Range("M138").Value = -1083.625100000001
$ws.Range("N138").Value = -22772.6885

# --- Sheet: ARM (43 cell updates) ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 9550.35
$ws.Range("I61").Value = 4577.522
$ws.Range("J61").Value = 16278.294
$ws.Range("K61").Value = 4577.522
$ws.Range("L61").Value = 16278.294
$ws.Range("M61").Value = -4365.522
$ws.Range("N61").Value = -16702.294
$ws.Range("H97").Value = 1069.375
$ws.Range("I97").Value = 824.36
$ws.Range("J97").Value = 1944.4286
$ws.Range("K97").Value = 824.36
$ws.Range("L97").Value = 1944.4286
$ws.Range("M97").Value = -328.36
$ws.Range("N97").Value = -2936.4286
$ws.Range("H122").Value = 11367445
$ws.Range("I122").Value = 4137.375
$ws.Range("J122").Value = 41669600
$ws.Range("K122").Value = 12412.125
$ws.Range("L122").Value = 125008800
$ws.Range("M122").Value = -9962.125
$ws.Range("N122").Value = -125013700
$ws.Range("H124").Value = 30000
$ws.Range("J124").Value = 30000
$ws.Range("L124").Value = 30000
$ws.Range("N124").Value = -39820
$ws.Range("H125").Value = 65436.363
$ws.Range("J125").Value = 65436.363
$ws.Range("L125").Value = 65436.363
$ws.Range("N125").Value = -75276.363
$ws.Range("H132").Value = 6272.3823
$ws.Range("I132").Value = 2455.8125
$ws.Range("J132").Value = 9664.888999999999
$ws.Range("K132").Value = 7367.4375
$ws.Range("L132").Value = 28994.667
$ws.Range("M132").Value = -4837.4375
$ws.Range("N132").Value = -34054.667
$ws.Range("H136").Value = 9550.35
$ws.Range("I136").Value = 4577.522
$ws.Range("J136").Value = 16278.294
$ws.Range("K136").Value = 13732.566
$ws.Range("L136").Value = 48834.882
$ws.Range("M136").Value = -11182.566
$ws.Range("N136").Value = -53934.882

# --- Sheet: BSM (15 cell updates) ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H15").Value = 6999.75
$ws.Range("J15").Value = 6999.75
$ws.Range("L15").Value = 6999.75
$ws.Range("N15").Value = -7453.75
$ws.Range("H19").Value = 20000
$ws.Range("J19").Value = 20000
$ws.Range("L19").Value = 20000
$ws.Range("N19").Value = -20346
$ws.Range("H134").Value = 30522.084
$ws.Range("I134").Value = 2711.7812
$ws.Range("J134").Value = 253004.5
$ws.Range("K134").Value = 8135.3436
$ws.Range("L134").Value = 759013.5
$ws.Range("M134").Value = -5600.3436
$ws.Range("N134").Value = -764083.5

# --- Sheet: CRP (26 cell updates) ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1717.28
$ws.Range("I31").Value = 1265.0632
$ws.Range("J31").Value = 3418.476
$ws.Range("K31").Value = 1265.0632
$ws.Range("L31").Value = 3418.476
$ws.Range("M31").Value = -970.0632000000001
$ws.Range("N31").Value = -4008.476
$ws.Range("H34").Value = 1717.28
$ws.Range("I34").Value = 1265.0632
$ws.Range("J34").Value = 3418.476
$ws.Range("K34").Value = 1265.0632
$ws.Range("L34").Value = 3418.476
$ws.Range("M34").Value = -1063.0632
$ws.Range("N34").Value = -3822.476
$ws.Range("H70").Value = 0
$ws.Range("J70").Value = 0
$ws.Range("L70").Value = 0
$ws.Range("N70").ClearContents()
$ws.Range("H73").Value = 0
$ws.Range("J73").Value = 0
$ws.Range("L73").Value = 0
$ws.Range("N73").ClearContents()
$ws.Range("H102").Value = 0
$ws.Range("J102").Value = 0
$ws.Range("L102").Value = 0
$ws.Range("N102").ClearContents()

# --- Sheet: CUL (50 cell updates) ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 1063.9286
$ws.Range("I23").Value = 4014
$ws.Range("J23").Value = 259.36365
$ws.Range("K23").Value = 12042
$ws.Range("L23").Value = 778.09095
$ws.Range("M23").Value = -11807
$ws.Range("N23").Value = -1248.09095
$ws.Range("H64").Value = 3225.3333
$ws.Range("I64").Value = 1730
$ws.Range("J64").Value = 3723.7778
$ws.Range("K64").Value = 5190
$ws.Range("L64").Value = 11171.3334
$ws.Range("M64").Value = -4920
$ws.Range("N64").Value = -11711.3334
$ws.Range("H67").Value = 3225.3333
$ws.Range("I67").Value = 1730
$ws.Range("J67").Value = 3723.7778
$ws.Range("K67").Value = 5190
$ws.Range("L67").Value = 11171.3334
$ws.Range("M67").Value = -4254
$ws.Range("N67").Value = -13043.3334
$ws.Range("H87").Value = 7830.1665
$ws.Range("I87").Value = 1333.3334
$ws.Range("J87").Value = 9995.777
$ws.Range("K87").Value = 4000.0002
$ws.Range("L87").Value = 29987.331
$ws.Range("M87").Value = -2752.0002
$ws.Range("N87").Value = -32483.331
$ws.Range("H90").Value = 7830.1665
$ws.Range("I90").Value = 1333.3334
$ws.Range("J90").Value = 9995.777
$ws.Range("K90").Value = 12000.0006
$ws.Range("L90").Value = 89961.993
$ws.Range("M90").Value = -5760.000599999999
$ws.Range("N90").Value = -102441.993
$ws.Range("H122").Value = 904.5238000000001
$ws.Range("J122").Value = 1399.8
$ws.Range("L122").Value = 12598.2
$ws.Range("N122").Value = -17498.2
$ws.Range("H132").Value = 2034.9231
$ws.Range("I132").Value = 2282
$ws.Range("J132").Value = 1746.6666
$ws.Range("K132").Value = 20538
$ws.Range("L132").Value = 15719.9994
$ws.Range("M132").Value = -18008
$ws.Range("N132").Value = -20779.9994
$ws.Range("H138").Value = 4879.9443
$ws.Range("I138").Value = 8187.143
$ws.Range("K138").Value = 24561.429
$ws.Range("M138").Value = -19421.429

# --- Sheet: GSM (12 cell updates) ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H27").Value = 9280
$ws.Range("J27").Value = 9280
$ws.Range("L27").Value = 9280
$ws.Range("N27").Value = -9612
$ws.Range("H104").Value = 59300
$ws.Range("J104").Value = 59300
$ws.Range("L104").Value = 59300
$ws.Range("N104").Value = -66288
$ws.Range("H122").Value = 7334.778
$ws.Range("I122").Value = 7916.143
$ws.Range("K122").Value = 23748.429
$ws.Range("M122").Value = -21298.429

# --- Sheet: LTW (7 cell updates) ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 308281.78
$ws.Range("I55").Value = 500758.38
$ws.Range("J55").Value = 319.2
$ws.Range("K55").Value = 500758.38
$ws.Range("L55").Value = 319.2
$ws.Range("M55").Value = -500585.38
$ws.Range("N55").Value = -665.2

# --- Sheet: WVR (7 cell updates) ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 4459.3066
$ws.Range("I136").Value = 1792.8918
$ws.Range("J136").Value = 8405.6
$ws.Range("K136").Value = 5378.6754
$ws.Range("L136").Value = 25216.8
$ws.Range("M136").Value = -2828.6754
$ws.Range("N136").Value = -30316.8
